$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 32 (shifts rows 32..332 down to 33..333),
# matching the "strWindowPos" entry added to the localization table.
$ws.Rows.Item(32).Insert()

# Populate the freshly inserted row 32 (order chosen so new shared
# strings are appended in the same order the target workbook uses).
$ws.Range("B32").Value = "localization\strings"
$ws.Range("C32").Value = "strWindowPos"
$ws.Range("D32").Value = "In ""settings"" form, tab ""User interface"""
$ws.Range("E32").Value = "Remember window position and size on startup"

# Row 25's Comment column gets the same new string as D32.
$ws.Range("D25").Value = "In ""settings"" form, tab ""User interface"""

# Grow the table (ListObject) so the new row is included.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("B2:F204"))

# Column D widened (auto-fit side effect of the new, longer comment text).
$ws.Columns.Item(4).ColumnWidth = 34.85
